$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value (45178) for every data
# row (rows 2 through 232). Bump every one of those cells to 45179.
$range = $ws.Range("C2:C232")
$range.Value = 45179
